# Finish To Add King to All AI functions
# Marks several more AI functions as "done" (strike-through), marks
# AISimpleMove as "in progress" (underline), and moves the Word
# "_GoBack" bookmark from the old "Drop" entry to the new last-edited
# entry "AI_eat_again".

$d = $word.ActiveDocument

# Locate paragraphs by their text so the script is resilient to any
# re-numbering.
function Get-ParaByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd("`r")
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# Move the hidden "_GoBack" bookmark off of "Drop" - it will be
# re-created further down, anchored on "AI_eat_again".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# AICheckEatFirst -> strike-through (done)
$pAICheckEatFirst = Get-ParaByText $d "AICheckEatFirst"
$pAICheckEatFirst.Range.Font.StrikeThrough = 1

# AISimpleMove -> underline (in progress)
$pAISimpleMove = Get-ParaByText $d "AISimpleMove"
$pAISimpleMove.Range.Font.Underline = 1

# AI_turn_start -> strike-through (done)
$pAITurnStart = Get-ParaByText $d "AI_turn_start"
$pAITurnStart.Range.Font.StrikeThrough = 1

# AI_eat_again -> strike-through (done) + re-anchor "_GoBack" here
$pAIEatAgain = Get-ParaByText $d "AI_eat_again"
$pAIEatAgain.Range.Font.StrikeThrough = 1
$d.Bookmarks.Add("_GoBack", $pAIEatAgain.Range)

# AI_Eat_Move -> strike-through (done)
$pAIEatMove = Get-ParaByText $d "AI_Eat_Move"
$pAIEatMove.Range.Font.StrikeThrough = 1

# AI_Move -> strike-through (done)
$pAIMove = Get-ParaByText $d "AI_Move"
$pAIMove.Range.Font.StrikeThrough = 1
